# Generate Report for Handback
# - Updates the existing handback entry for 49e75141-5af2-4885-aca5-2eac0a7b72b5.md
#   to now point at acb3d08b-601e-4505-b3a3-5b94ba208151.md (new hashes / timestamps).
# - Adds a brand new handback entry for ae20c40a-579e-4708-88c6-9d041cfce420.md
#   as a second row on every sheet (Overview, zh-cn, de-de).

function Set-Text($ws, $addr, $val) {
    # Leading apostrophe forces literal text entry (prevents Excel's automatic
    # Boolean / number / empty-string coercion from kicking in).
    $ws.Range($addr).Value = "'" + $val
}

function Set-DateText($ws, $addr, $val) {
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Existing row 2: rename the referenced file.
Set-Text   $wsOverview "A2" "acb3d08b-601e-4505-b3a3-5b94ba208151.md"
Set-Text   $wsOverview "B2" "e2e\acb3d08b-601e-4505-b3a3-5b94ba208151.md"
Set-DateText $wsOverview "G2" "2016-09-05 05:07:49"

# New row 3 for the second handed-back file.
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

Set-Text   $wsOverview "A3" "ae20c40a-579e-4708-88c6-9d041cfce420.md"
Set-Text   $wsOverview "B3" "e2e\ae20c40a-579e-4708-88c6-9d041cfce420.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/ae20c40a-579e-4708-88c6-9d041cfce420.md", "", "", "e2e\ae20c40a-579e-4708-88c6-9d041cfce420.md") | Out-Null
Set-Text   $wsOverview "C3" ".md"
Set-Text   $wsOverview "E3" "Handed back: in sync with en-US"
Set-Text   $wsOverview "F3" "Handed back: in sync with en-US"
Set-DateText $wsOverview "G3" "2016-09-05 05:07:49"

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Existing row 2: rename the referenced file + refresh its handoff/handback timestamps.
Set-Text     $wsZh "A2" "acb3d08b-601e-4505-b3a3-5b94ba208151.md"
Set-Text     $wsZh "G2" "acb3d08b-601e-4505-b3a3-5b94ba208151.dafc2fc1a905c76ab6110551a560695cdf4f3527.zh-cn.xlf"
Set-DateText $wsZh "H2" "2016-09-05 05:07:44"
Set-Text     $wsZh "I2" "acb3d08b-601e-4505-b3a3-5b94ba208151.md"
Set-Text     $wsZh "J2" "acb3d08b-601e-4505-b3a3-5b94ba208151.dafc2fc1a905c76ab6110551a560695cdf4f3527.zh-cn.xlf"
Set-DateText $wsZh "K2" "2016-09-05 05:08:06"

# New row 3 for the second handed-back file.
$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null

Set-Text $wsZh "A3" "ae20c40a-579e-4708-88c6-9d041cfce420.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/ae20c40a-579e-4708-88c6-9d041cfce420.md", "", "", "ae20c40a-579e-4708-88c6-9d041cfce420.md") | Out-Null
Set-Text     $wsZh "B3" ".md"
Set-Text     $wsZh "C3" "Handed back: in sync with en-US"
Set-Text     $wsZh "D3" "e2e"
Set-Text     $wsZh "E3" "ht"
Set-Text     $wsZh "F3" "True"
Set-Text     $wsZh "G3" "ae20c40a-579e-4708-88c6-9d041cfce420.c9cd2fe18d1aaefff2835fed38c90737da6cd0fd.zh-cn.xlf"
Set-DateText $wsZh "H3" "2016-09-05 05:07:44"
Set-Text $wsZh "I3" "ae20c40a-579e-4708-88c6-9d041cfce420.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/master/e2e/ae20c40a-579e-4708-88c6-9d041cfce420.md", "", "", "ae20c40a-579e-4708-88c6-9d041cfce420.md") | Out-Null
Set-Text     $wsZh "J3" "ae20c40a-579e-4708-88c6-9d041cfce420.c9cd2fe18d1aaefff2835fed38c90737da6cd0fd.zh-cn.xlf"
Set-DateText $wsZh "K3" "2016-09-05 05:08:06"
Set-Text     $wsZh "L3" ""
Set-Text     $wsZh "M3" "True"
Set-Text     $wsZh "N3" ""
Set-Text     $wsZh "O3" "False"
Set-Text     $wsZh "P3" ""

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Existing row 2: rename the referenced file + refresh its handoff/handback timestamps.
Set-Text     $wsDe "A2" "acb3d08b-601e-4505-b3a3-5b94ba208151.md"
Set-Text     $wsDe "G2" "acb3d08b-601e-4505-b3a3-5b94ba208151.dafc2fc1a905c76ab6110551a560695cdf4f3527.de-de.xlf"
Set-Text     $wsDe "I2" "acb3d08b-601e-4505-b3a3-5b94ba208151.md"
Set-Text     $wsDe "J2" "acb3d08b-601e-4505-b3a3-5b94ba208151.dafc2fc1a905c76ab6110551a560695cdf4f3527.de-de.xlf"
Set-DateText $wsDe "K2" "2016-09-05 05:08:17"

# New row 3 for the second handed-back file.
$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null

Set-Text $wsDe "A3" "ae20c40a-579e-4708-88c6-9d041cfce420.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/ae20c40a-579e-4708-88c6-9d041cfce420.md", "", "", "ae20c40a-579e-4708-88c6-9d041cfce420.md") | Out-Null
Set-Text     $wsDe "B3" ".md"
Set-Text     $wsDe "C3" "Handed back: in sync with en-US"
Set-Text     $wsDe "D3" "e2e"
Set-Text     $wsDe "E3" "ht"
Set-Text     $wsDe "F3" "True"
Set-Text     $wsDe "G3" "ae20c40a-579e-4708-88c6-9d041cfce420.c9cd2fe18d1aaefff2835fed38c90737da6cd0fd.de-de.xlf"
Set-DateText $wsDe "H3" "2016-09-05 05:07:49"
Set-Text $wsDe "I3" "ae20c40a-579e-4708-88c6-9d041cfce420.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/master/e2e/ae20c40a-579e-4708-88c6-9d041cfce420.md", "", "", "ae20c40a-579e-4708-88c6-9d041cfce420.md") | Out-Null
Set-Text     $wsDe "J3" "ae20c40a-579e-4708-88c6-9d041cfce420.c9cd2fe18d1aaefff2835fed38c90737da6cd0fd.de-de.xlf"
Set-DateText $wsDe "K3" "2016-09-05 05:08:17"
Set-Text     $wsDe "L3" ""
Set-Text     $wsDe "M3" "True"
Set-Text     $wsDe "N3" ""
Set-Text     $wsDe "O3" "False"
Set-Text     $wsDe "P3" ""
